$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row based on the sheet's dimension (column C holds the "Förändrad" date)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
